# Updated Argent prices: append a new "2025-04-28" row to every price
# sheet, carrying forward the previous day's (last row's) price value,
# exactly as in the source diff (new row 58, dimension A1:B57 -> A1:B58).

$wb = $excel.ActiveWorkbook
$newDate = "2025-04-28"

foreach ($ws in $wb.Worksheets) {
    # Column A holds the dates as text; find the last used row by walking
    # up from the bottom of the sheet (mirrors Ctrl+Up in column A).
    $lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
    $newRow = $lastRow + 1

    # Read the previous day's price as displayed text (.Value on these
    # cells can misbehave in this host, .Text is reliable) and carry it
    # forward unchanged into the new row.
    $priorPrice = $ws.Cells($lastRow, 2).Text

    $dateCell = $ws.Cells($newRow, 1)
    $priceCell = $ws.Cells($newRow, 2)

    # Force text storage (matching the rest of the column, which holds
    # literal text rather than real dates/numbers) by pre-formatting the
    # cells as Text before writing, then clearing the format back off so
    # the new cells don't end up with a stray number-format style applied
    # - mirrors the surrounding cells, which carry no explicit style.
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = $newDate
    $priceCell.Value = $priorPrice

    $ws.Range($dateCell, $priceCell).ClearFormats()
}
